# Scheduled-runner data refresh: update Universalis market-price derived
# columns (H:N) on each job sheet. Values only -- no formulas, styles, or
# structural changes involved.
$wb = $excel.ActiveWorkbook

# --- ALC ---
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H76").Value = 3195
$ws.Range("I76").Value = 3195
$ws.Range("K76").Value = 3195
$ws.Range("M76").Value = -2880
$ws.Range("H79").Value = 3195
$ws.Range("I79").Value = 3195
$ws.Range("K79").Value = 3195
$ws.Range("M79").Value = -2103
$ws.Range("H98").Value = 2483.121
$ws.Range("I98").Value = 2593.261
$ws.Range("J98").Value = 2229.8
$ws.Range("K98").Value = 2593.261
$ws.Range("L98").Value = 2229.8
$ws.Range("M98").Value = -1095.261
$ws.Range("N98").Value = -5225.8
$ws.Range("H122").Value = 2483.121
$ws.Range("I122").Value = 2593.261
$ws.Range("J122").Value = 2229.8
$ws.Range("K122").Value = 7779.782999999999
$ws.Range("L122").Value = 6689.400000000001
$ws.Range("M122").Value = -5329.782999999999
$ws.Range("N122").Value = -11589.4
$ws.Range("H127").Value = 1494.95
$ws.Range("I127").Value = 1099.9
$ws.Range("J127").Value = 1890
$ws.Range("K127").Value = 3299.7
$ws.Range("L127").Value = 5670
$ws.Range("M127").Value = 1660.3
$ws.Range("N127").Value = -15590
$ws.Range("H129").Value = 1423.7941
$ws.Range("I129").Value = 740.6667
$ws.Range("J129").Value = 1570.1786
$ws.Range("K129").Value = 2222.0001
$ws.Range("L129").Value = 4710.5358
$ws.Range("M129").Value = 2777.9999
$ws.Range("N129").Value = -14710.5358
$ws.Range("H131").Value = 1643.7778
$ws.Range("I131").Value = 1699.25
$ws.Range("J131").Value = 1200
$ws.Range("K131").Value = 5097.75
$ws.Range("L131").Value = 3600
$ws.Range("M131").Value = -57.75
$ws.Range("N131").Value = -13680
$ws.Range("H132").Value = 273584.75
$ws.Range("I132").Value = 281156.38
$ws.Range("J132").Value = 1006
$ws.Range("K132").Value = 843469.14
$ws.Range("L132").Value = 3018
$ws.Range("M132").Value = -840939.14
$ws.Range("N132").Value = -8078
$ws.Range("H135").Value = 1955.8667
$ws.Range("I135").Value = 690.7083
$ws.Range("J135").Value = 7016.5
$ws.Range("K135").Value = 6216.3747
$ws.Range("L135").Value = 63148.5
$ws.Range("M135").Value = -3681.3747
$ws.Range("N135").Value = -68218.5
$ws.Range("H137").Value = 1669.1147
$ws.Range("I137").Value = 1234.7
$ws.Range("J137").Value = 2089.516
$ws.Range("K137").Value = 3704.1
$ws.Range("L137").Value = 6268.548000000001
$ws.Range("M137").Value = -1154.1
$ws.Range("N137").Value = -11368.548

# --- ARM ---
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H32").Value = 11394.121
$ws.Range("I32").Value = 4950.241
$ws.Range("J32").Value = 58112.25
$ws.Range("K32").Value = 4950.241
$ws.Range("L32").Value = 58112.25
$ws.Range("M32").Value = -4663.241
$ws.Range("N32").Value = -58686.25
$ws.Range("H61").Value = 2188.25
$ws.Range("I61").Value = 1749.6857
$ws.Range("J61").Value = 3369
$ws.Range("K61").Value = 1749.6857
$ws.Range("L61").Value = 3369
$ws.Range("M61").Value = -1537.6857
$ws.Range("N61").Value = -3793
$ws.Range("H74").Value = 5157.913
$ws.Range("I74").Value = 756.13336
$ws.Range("J74").Value = 13411.25
$ws.Range("K74").Value = 756.13336
$ws.Range("L74").Value = 13411.25
$ws.Range("M74").Value = 117.86664
$ws.Range("N74").Value = -15159.25
$ws.Range("H77").Value = 5157.913
$ws.Range("I77").Value = 756.13336
$ws.Range("J77").Value = 13411.25
$ws.Range("K77").Value = 3780.6668
$ws.Range("L77").Value = 67056.25
$ws.Range("M77").Value = 587.3332
$ws.Range("N77").Value = -75792.25
$ws.Range("H136").Value = 2188.25
$ws.Range("I136").Value = 1749.6857
$ws.Range("J136").Value = 3369
$ws.Range("K136").Value = 5249.0571
$ws.Range("L136").Value = 10107
$ws.Range("M136").Value = -2699.0571
$ws.Range("N136").Value = -15207

# --- BSM ---
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H94").Value = 1009.9355
$ws.Range("I94").Value = 969.0952
$ws.Range("K94").Value = 969.0952
$ws.Range("M94").Value = -518.0952
$ws.Range("H134").Value = 6243
$ws.Range("I134").Value = 7003.9473
$ws.Range("J134").Value = 2628.5
$ws.Range("K134").Value = 21011.8419
$ws.Range("L134").Value = 7885.5
$ws.Range("M134").Value = -18476.8419
$ws.Range("N134").Value = -12955.5

# --- CRP ---
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H31").Value = 1446.1482
$ws.Range("I31").Value = 1132.5483
$ws.Range("J31").Value = 1868.826
$ws.Range("K31").Value = 1132.5483
$ws.Range("L31").Value = 1868.826
$ws.Range("M31").Value = -837.5482999999999
$ws.Range("N31").Value = -2458.826
$ws.Range("H34").Value = 1446.1482
$ws.Range("I34").Value = 1132.5483
$ws.Range("J34").Value = 1868.826
$ws.Range("K34").Value = 1132.5483
$ws.Range("L34").Value = 1868.826
$ws.Range("M34").Value = -930.5482999999999
$ws.Range("N34").Value = -2272.826
$ws.Range("H132").Value = 3172.8276
$ws.Range("I132").Value = 2839.2856
$ws.Range("J132").Value = 4048.375
$ws.Range("K132").Value = 8517.856800000001
$ws.Range("L132").Value = 12145.125
$ws.Range("M132").Value = -5987.856800000001
$ws.Range("N132").Value = -17205.125

# --- CUL ---
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H4").Value = 1167.6666
$ws.Range("I4").Value = 1001
$ws.Range("J4").Value = 1251
$ws.Range("K4").Value = 3003
$ws.Range("L4").Value = 3753
$ws.Range("M4").Value = -2891
$ws.Range("N4").Value = -3977
$ws.Range("H68").Value = 950.5679
$ws.Range("I68").Value = 726.29266
$ws.Range("J68").Value = 1180.45
$ws.Range("K68").Value = 2178.87798
$ws.Range("L68").Value = 3541.35
$ws.Range("M68").Value = -1367.87798
$ws.Range("N68").Value = -5163.35
$ws.Range("H71").Value = 950.5679
$ws.Range("I71").Value = 726.29266
$ws.Range("J71").Value = 1180.45
$ws.Range("K71").Value = 6536.63394
$ws.Range("L71").Value = 10624.05
$ws.Range("M71").Value = -2480.63394
$ws.Range("N71").Value = -18736.05
$ws.Range("H107").Value = 54769.49
$ws.Range("I107").Value = 31889.656
$ws.Range("J107").Value = 201200.4
$ws.Range("K107").Value = 95668.96799999999
$ws.Range("L107").Value = 603601.2
$ws.Range("M107").Value = -93748.96799999999
$ws.Range("N107").Value = -607441.2
$ws.Range("H131").Value = 1353145.4
$ws.Range("I131").Value = 1115.6428
$ws.Range("J131").Value = 1668618.9
$ws.Range("K131").Value = 3346.9284
$ws.Range("L131").Value = 5005856.699999999
$ws.Range("M131").Value = 1693.0716
$ws.Range("N131").Value = -5015936.699999999

# --- GSM ---
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H49").Value = 5000
$ws.Range("J49").Value = 5000
$ws.Range("L49").Value = 5000
$ws.Range("N49").Value = -5368
$ws.Range("H122").Value = 6563.6
$ws.Range("I122").Value = 7192
$ws.Range("J122").Value = 4050
$ws.Range("K122").Value = 21576
$ws.Range("L122").Value = 12150
$ws.Range("M122").Value = -19126
$ws.Range("N122").Value = -17050
$ws.Range("H132").Value = 2653.5
$ws.Range("I132").Value = 3095.5334
$ws.Range("J132").Value = 1706.2858
$ws.Range("K132").Value = 9286.600199999999
$ws.Range("L132").Value = 5118.857400000001
$ws.Range("M132").Value = -6756.600199999999
$ws.Range("N132").Value = -10178.8574

# --- LTW ---
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H40").Value = 3182.3333
$ws.Range("I40").Value = 2773.5
$ws.Range("K40").Value = 2773.5
$ws.Range("M40").Value = -2637.5
$ws.Range("H132").Value = 8115.9
$ws.Range("I132").Value = 11745.818
$ws.Range("J132").Value = 3679.3333
$ws.Range("K132").Value = 35237.454
$ws.Range("L132").Value = 11037.9999
$ws.Range("M132").Value = -32707.454
$ws.Range("N132").Value = -16097.9999

# --- WVR ---
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H132").Value = 4547.8335
$ws.Range("I132").Value = 4944.696
$ws.Range("J132").Value = 3243.8572
$ws.Range("K132").Value = 14834.088
$ws.Range("L132").Value = 9731.571599999999
$ws.Range("M132").Value = -12304.088
$ws.Range("N132").Value = -14791.5716
$ws.Range("H136").Value = 853.7879
$ws.Range("I136").Value = 750.55554
$ws.Range("J136").Value = 1318.3334
$ws.Range("K136").Value = 2251.66662
$ws.Range("L136").Value = 3955.0002
$ws.Range("M136").Value = 298.33338
$ws.Range("N136").Value = -9055.0002
